$wb = $excel.ActiveWorkbook

# Both "展览" and "全部类型" sheets contain identical data tables that
# need their "想去人数" (F column) counts refreshed.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F9").Value = 6097
    $ws.Range("F16").Value = 775
    $ws.Range("F17").Value = 175
}
